$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Permits Filed for 147-02 106th Avenue in Jamaica, Queens"
$ws.Range("B2").Value = "https://newyorkyimby.com/2026/02/permits-filed-for-147-02-106th-avenue-in-jamaica-queens.html"
$ws.Range("C2").Value = "Permits have been filed for an 11-story residential building at 147-02 106th Avenue in <a href=""https://newyorkyimby.com/neighborhoods/jamaica"">Jamaica</a>, Queens. Also addressed as 106-03 Sutphin Boulevard, the corner lot is near the Sutphin Boulevard–Archer Avenue–JFK Airport subway station, served by the E, J, and Z trains. Xue Mei Yi of 12001 Realty LLC is listed as the owner behind the applications."
$ws.Range("D2").Value = "2026-02-11T11:30:20+00:00"
$ws.Range("E2").Value = "Wed, 11 Feb 2026 11:30:20 +0000"
$ws.Range("F2").Value = "YIMBY"
$ws.Range("G2").Value = "YIMBY - Jamaica"
$ws.Range("H2").Value = ""
